$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates (Fgf10-Fgfr1 -> ECs TPM values)
$ws.Range("M2").Value = 5.978421000000001
$ws.Range("N2").Value = 17.935263
$ws.Range("O2").Value = 0.05704457007880161
$ws.Range("P2").Value = 0.06242884486533885
$ws.Range("Q2").Value = 1.58850631584
$ws.Range("R2").Value = 14.29655684256
$ws.Range("S2").Value = 0.05704457007880161
$ws.Range("T2").Value = 0.06242884486533885

# Row 3 updates (Fgf10-Fgfr1 -> FAPs TPM values)
$ws.Range("O3").Value = 0.6646576013185088
$ws.Range("P3").Value = 0.7273927426214574
$ws.Range("S3").Value = 0.6646576013185088
$ws.Range("T3").Value = 0.7273927426214574

# Row 4 updates (Fgf10-Fgfr1 -> Inflammatory-Mac TPM values)
$ws.Range("M4").Value = 1.290243
$ws.Range("N4").Value = 3.870729
$ws.Range("O4").Value = 0.01231116999491725
$ws.Range("P4").Value = 0.01347318632889677
$ws.Range("Q4").Value = 0.34282616672
$ws.Range("R4").Value = 3.08543550048
$ws.Range("S4").Value = 0.01231116999491725
$ws.Range("T4").Value = 0.01347318632889677

# Row 5 updates (Fgf10-Fgfr1 -> MuSCs TPM values)
$ws.Range("M5").Value = 27.1166075
$ws.Range("N5").Value = 54.233215
$ws.Range("O5").Value = 0.2587397603536297
$ws.Range("P5").Value = 0.1887743138075849
$ws.Range("Q5").Value = 7.205063390133333
$ws.Range("R5").Value = 43.2303803408
$ws.Range("S5").Value = 0.2587397603536297
$ws.Range("T5").Value = 0.1887743138075849

# Row 6 updates (Fgf10-Fgfr1 -> Resolving-Mac TPM values)
$ws.Range("M6").Value = 0.759494
$ws.Range("N6").Value = 2.278482
$ws.Range("O6").Value = 0.00724689825414258
$ws.Range("P6").Value = 0.007930912376722157
$ws.Range("Q6").Value = 0.2018026190933333
$ws.Range("R6").Value = 1.81622357184
$ws.Range("S6").Value = 0.00724689825414258
$ws.Range("T6").Value = 0.007930912376722157
